$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.113.45'
$ws.Range("E2").Value = '  -4.01%  '

$ws.Range("D3").Value = '2.620.13'
$ws.Range("E3").Value = '  -3.13%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").Value = "'517.50"
$ws.Range("E5").Value = '  -1.84%  '

$ws.Range("D6").Value = "'141.77"
$ws.Range("E6").Value = '  -2.15%  '

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = '  +0.29%  '

$ws.Range("E8").Value = '  -2.11%  '

$ws.Range("D9").Value = "'6.66"
$ws.Range("E9").Value = '  -0.64%  '

$ws.Range("E10").Value = '  -3.15%  '

$ws.Range("D11").Value = "'0.335"
$ws.Range("E11").Value = '  -1.40%  '

$ws.Range("E12").Value = '  +1.28%  '

$ws.Range("D13").Value = '3.078.93'
$ws.Range("E13").Value = '  -3.32%  '

$ws.Range("D14").Value = '58.075.21'
$ws.Range("E14").Value = '  -4.10%  '

$ws.Range("D15").Value = "'20.66"
$ws.Range("E15").Value = '  -3.14%  '

$ws.Range("E16").Value = '  -1.74%  '

$ws.Range("D17").Value = '2.621.32'
$ws.Range("E17").Value = '  -8.42%  '

$ws.Range("E18").Value = '  -2.93%  '

$ws.Range("D19").Value = "'334.61"
$ws.Range("E19").Value = '  -3.51%  '

$ws.Range("D20").Value = "'10.35"
$ws.Range("E20").Value = '  -2.82%  '

$ws.Range("D21").Value = "'6.26"
$ws.Range("E21").Value = '  -2.95%  '

$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = '  +0.15%  '

$ws.Range("D23").Value = "'64.16"
$ws.Range("E23").Value = '  +0.71%  '

$ws.Range("D24").Value = "'0.423"
$ws.Range("E24").Value = '  +0.59%  '

$ws.Range("E25").Value = '  -2.43%  '

$ws.Range("E26").Value = '  +0.62%  '

$ws.Range("D27").Value = "'7.05"
$ws.Range("E27").Value = '  -3.29%  '

$ws.Range("D28").Value = '0.0₃0783'

$ws.Range("D29").Value = "'6.60"
$ws.Range("E29").Value = '  -2.89%  '

$ws.Range("E30").Value = '  +0.08%  '

$ws.Range("D31").Value = "'152.38"
$ws.Range("E31").Value = '  +1.51%  '

$ws.Range("E32").Value = '  -1.51%  '

$ws.Range("D33").Value = "'18.70"
$ws.Range("E33").Value = '  -2.06%  '

$ws.Range("D34").Value = "'4.09"
$ws.Range("E34").Value = '  -3.80%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = "'1.17"
$ws.Range("E35").Value = '  -4.93%  '

$ws.Range("B36").Value = 'SuiNetwork'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D36").Value = "'0.899"
$ws.Range("E36").Value = '  -4.54%  '

$ws.Range("D37").Value = "'36.65"
$ws.Range("E37").Value = '  -1.23%  '

$ws.Range("D38").Value = "'0.846"
$ws.Range("E38").Value = '  -3.23%  '

$ws.Range("D39").Value = "'1.43"
$ws.Range("E39").Value = '  -5.69%  '

$ws.Range("D40").Value = "'3.60"
$ws.Range("E40").Value = '  -1.84%  '

$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = '  +0.38%  '

$ws.Range("D42").Value = "'0.598"
$ws.Range("E42").Value = '  -2.01%  '

$ws.Range("D43").Value = "'0.0966"
$ws.Range("E43").Value = '  -2.40%  '

$ws.Range("D44").Value = "'268.39"
$ws.Range("E44").Value = '  -5.37%  '

$ws.Range("B45").Value = 'WhiteBITCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D45").Value = "'10.61"
$ws.Range("E45").Value = '  +1.30%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = "'19.20"
$ws.Range("E46").Value = '  -4.53%  '

$ws.Range("E47").Value = '  -1.47%  '

$ws.Range("D48").Value = '2.031.32'
$ws.Range("E48").Value = '  -5.19%  '

$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = "'0.0226"
$ws.Range("E49").Value = '  -3.13%  '

$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = "'4.64"
$ws.Range("E50").Value = '  -4.11%  '

$ws.Range("D51").Value = "'18.22"
$ws.Range("E51").Value = '  -5.11%  '
